$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values for the crypto symbol list snapshot.
# Cells are stored as text in the workbook (Price / Volume columns),
# so we force text using a leading apostrophe and then reset the
# cell style back to "Normal" so no stray number-format / quote-prefix
# style is left behind on the cell.
$updates = @{
    'D2' = '330.42'
    'E2' = '0.64%'
    'D4' = '5.721'
    'E4' = '-1.85%'
    'D5' = '0.08105'
    'E5' = '0.89%'
    'D6' = '2.040'
    'E6' = '5.48%'
    'D7' = '8.740'
    'E8' = '-1.37%'
    'D10' = '0.9233'
    'E10' = '-2.25%'
    'E11' = '0.21%'
    'D12' = '0.1946'
    'E12' = '-0.79%'
    'D13' = '8.328'
    'E13' = '-6.26%'
    'D14' = '0.09281'
    'E14' = '0.89%'
    'D15' = '0.03663'
    'E15' = '1.96%'
    'D16' = '0.1056'
    'E16' = '9.62%'
    'D17' = '0.001299'
    'E17' = '-2.52%'
    'D18' = '0.006250'
    'E18' = '-1.54%'
    'D19' = '3.384'
    'E19' = '0.41%'
    'E20' = '-1.23%'
    'D21' = '0.1416'
    'E21' = '-1.29%'
    'D22' = '0.2650'
    'E22' = '9.70%'
    'D23' = '0.04428'
    'E23' = '0.76%'
    'D24' = '0.001260'
    'E24' = '0.00%'
    'D25' = '0.004299'
    'E25' = '-0.87%'
    'E26' = '8.66%'
    'D39' = '0.02816'
    'E39' = '16.14%'
    'D40' = '0.05475'
    'E40' = '3.32%'
    'D41' = '0.007602'
    'E41' = '1.56%'
    'D42' = '0.009944'
    'E42' = '12.45%'
    'D43' = '0.1423'
    'E43' = '0.09%'
    'D44' = '0.002120'
    'E44' = '0.81%'
    'E45' = '11.44%'
    'D46' = '0.00006754'
    'E46' = '-2.19%'
    'E47' = '-0.38%'
    'D48' = '0.002955'
    'E48' = '-6.40%'
    'D49' = '0.002279'
    'E49' = '59.88%'
    'D50' = '0.00002100'
    'E50' = '-0.38%'
    'D51' = '0.0002000'
    'E51' = '-0.38%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.Value = "'" + $updates[$cell]
    $range.Style = "Normal"
}
